$wb = $excel.ActiveWorkbook

# Sheet ALC (sheet1) row 34
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 10316.667
$ws.Range("I34").Value = 1466.6666
$ws.Range("J34").Value = 19166.666
$ws.Range("K34").Value = 1466.6666
$ws.Range("L34").Value = 19166.666
$ws.Range("M34").Value = -1263.6666
$ws.Range("N34").Value = -19572.666

# Sheet ALC (sheet1) row 36
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H36").Value = 10316.667
$ws.Range("I36").Value = 1466.6666
$ws.Range("J36").Value = 19166.666
$ws.Range("K36").Value = 1466.6666
$ws.Range("L36").Value = 19166.666
$ws.Range("M36").Value = -751.6666
$ws.Range("N36").Value = -20596.666

# Sheet ALC (sheet1) row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 71432024
$ws.Range("I64").Value = 166668670
$ws.Range("J64").Value = 4545
$ws.Range("K64").Value = 166668670
$ws.Range("L64").Value = 4545
$ws.Range("M64").Value = -166668422
$ws.Range("N64").Value = -5041

# Sheet ALC (sheet1) row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 71432024
$ws.Range("I67").Value = 166668670
$ws.Range("J67").Value = 4545
$ws.Range("K67").Value = 166668670
$ws.Range("L67").Value = 4545
$ws.Range("M67").Value = -166667812
$ws.Range("N67").Value = -6261

# Sheet ALC (sheet1) row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 859.13635
$ws.Range("I129").Value = 764.8889
$ws.Range("J129").Value = 924.38464
$ws.Range("K129").Value = 2294.6667
$ws.Range("L129").Value = 2773.15392
$ws.Range("M129").Value = 2705.3333
$ws.Range("N129").Value = -12773.15392

# Sheet ALC (sheet1) row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2054.1208
$ws.Range("I132").Value = 1576.6082
$ws.Range("J132").Value = 4132.706
$ws.Range("K132").Value = 4729.8246
$ws.Range("L132").Value = 12398.118
$ws.Range("M132").Value = -2199.8246
$ws.Range("N132").Value = -17458.118

# Sheet ARM (sheet2) row 23
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 68433.71
$ws.Range("I23").Value = 77505.75
$ws.Range("J23").Value = 56337.668
$ws.Range("K23").Value = 77505.75
$ws.Range("L23").Value = 56337.668
$ws.Range("M23").Value = -77246.75
$ws.Range("N23").Value = -56855.668

# Sheet CRP (sheet4) row 87
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H87").Value = 48000
$ws.Range("J87").Value = 48000
$ws.Range("L87").Value = 48000
$ws.Range("N87").Value = -50372

# Sheet CRP (sheet4) row 90
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H90").Value = 48000
$ws.Range("J90").Value = 48000
$ws.Range("L90").Value = 144000
$ws.Range("N90").Value = -155856

# Sheet CUL (sheet5) row 76
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 4000
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 12000
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -12766

# Sheet CUL (sheet5) row 79
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H79").Value = 4000
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 12000
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -14652

# Sheet CUL (sheet5) row 82
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 8000
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

# Sheet CUL (sheet5) row 85
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H85").Value = 8000
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

# Sheet CUL (sheet5) row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 677.25
$ws.Range("I107").Value = 163.33333
$ws.Range("J107").Value = 985.6
$ws.Range("K107").Value = 489.99999
$ws.Range("L107").Value = 2956.8
$ws.Range("M107").Value = 1430.00001
$ws.Range("N107").Value = -6796.8

# Sheet GSM (sheet6) row 19
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 4500
$ws.Range("I19").Value = 4500
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 4500
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -4212
$ws.Range("N19").ClearContents()

# Sheet GSM (sheet6) row 69
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 20000.666
$ws.Range("J69").Value = 20000.666
$ws.Range("L69").Value = 20000.666
$ws.Range("N69").Value = -21498.666

# Sheet GSM (sheet6) row 72
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H72").Value = 20000.666
$ws.Range("J72").Value = 20000.666
$ws.Range("L72").Value = 60001.99800000001
$ws.Range("N72").Value = -67489.998

# Sheet GSM (sheet6) row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2786.7144
$ws.Range("I80").Value = 2333.75
$ws.Range("J80").Value = 3390.6667
$ws.Range("K80").Value = 2333.75
$ws.Range("L80").Value = 3390.6667
$ws.Range("M80").Value = -1335.75
$ws.Range("N80").Value = -5386.6667

# Sheet GSM (sheet6) row 82
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 10298
$ws.Range("I82").Value = 10298
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 10298
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -9915
$ws.Range("N82").ClearContents()

# Sheet GSM (sheet6) row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2786.7144
$ws.Range("I83").Value = 2333.75
$ws.Range("J83").Value = 3390.6667
$ws.Range("K83").Value = 11668.75
$ws.Range("L83").Value = 16953.3335
$ws.Range("M83").Value = -6676.75
$ws.Range("N83").Value = -26937.3335

# Sheet GSM (sheet6) row 85
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H85").Value = 10298
$ws.Range("I85").Value = 10298
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 10298
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -8972
$ws.Range("N85").ClearContents()

# Sheet GSM (sheet6) row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 834471.8
$ws.Range("I113").Value = 1250970.2
$ws.Range("J113").Value = 1475
$ws.Range("K113").Value = 1250970.2
$ws.Range("L113").Value = 1475
$ws.Range("M113").Value = -1248800.2
$ws.Range("N113").Value = -5815

# Sheet LTW (sheet7) row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5274.6665
$ws.Range("I22").Value = 433
$ws.Range("J22").Value = 6081.6113
$ws.Range("K22").Value = 433
$ws.Range("L22").Value = 6081.6113
$ws.Range("M22").Value = -138
$ws.Range("N22").Value = -6671.6113

# Sheet LTW (sheet7) row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 5274.6665
$ws.Range("I27").Value = 433
$ws.Range("J27").Value = 6081.6113
$ws.Range("K27").Value = 433
$ws.Range("L27").Value = 6081.6113
$ws.Range("M27").Value = -326
$ws.Range("N27").Value = -6295.6113

# Sheet LTW (sheet7) row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2296
$ws.Range("I68").Value = 2233.3333
$ws.Range("J68").Value = 2390
$ws.Range("K68").Value = 2233.3333
$ws.Range("L68").Value = 2390
$ws.Range("M68").Value = -1484.3333
$ws.Range("N68").Value = -3888

# Sheet LTW (sheet7) row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2296
$ws.Range("I71").Value = 2233.3333
$ws.Range("J71").Value = 2390
$ws.Range("K71").Value = 11166.6665
$ws.Range("L71").Value = 11950
$ws.Range("M71").Value = -7422.666499999999
$ws.Range("N71").Value = -19438

# Sheet LTW (sheet7) row 75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H75").Value = 41086.5
$ws.Range("J75").Value = 41086.5
$ws.Range("L75").Value = 41086.5
$ws.Range("N75").Value = -42958.5

# Sheet LTW (sheet7) row 78
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H78").Value = 41086.5
$ws.Range("J78").Value = 41086.5
$ws.Range("L78").Value = 123259.5
$ws.Range("N78").Value = -132619.5

# Sheet LTW (sheet7) row 81
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# Sheet LTW (sheet7) row 84
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# Sheet WVR (sheet8) row 43
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 10433.333
$ws.Range("I43").Value = 5000
$ws.Range("J43").Value = 13150
$ws.Range("K43").Value = 5000
$ws.Range("L43").Value = 13150
$ws.Range("M43").Value = -4851
$ws.Range("N43").Value = -13448
